$d = $word.ActiveDocument

# --- Paragraph 3 (Italy intro paragraph) text edits ---

# "building on " -> "exploiting "
$r = $d.Content
$r.Find.Execute("building on ", $true, $false, $false, $false, $false, $true, 1, $false, "exploiting ", 2)

# "As NATO navies go" -> "Within NATO"
$r = $d.Content
$r.Find.Execute("As NATO navies go", $true, $false, $false, $false, $false, $true, 1, $false, "Within NATO", 2)

# "many countries, Italy has had" -> "many countries however, Italy has had"
$r = $d.Content
$r.Find.Execute("many countries, Italy has had", $true, $false, $false, $false, $false, $true, 1, $false, "many countries however, Italy has had", 2)

# "forces and her navy in particular." -> "forces and on her navy in particular."
$r = $d.Content
$r.Find.Execute("forces and her navy in particular.", $true, $false, $false, $false, $false, $true, 1, $false, "forces and on her navy in particular.", 2)

# --- Paragraph 4 (FRY / Operation Sharp Guard paragraph) text edit ---
# "...rival factions in the FRY (Serbia, Croatia, Bosnia etc), and had..."
#   -> "...rival nationalist and ethnic factions in the FRY, which had..."
$r = $d.Content
$r.Find.Execute("rival factions in the FRY (Serbia, Croatia, Bosnia etc), and had", $true, $false, $false, $false, $false, $true, 1, $false, "rival nationalist and ethnic factions in the FRY, which had", 2)

# --- Move the _GoBack bookmark from the end of paragraph 2 (after the
#     Garibaldi picture) to its new location in paragraph 4, right after
#     "...FRY, which" and before " had the side benefit...". Re-adding a
#     bookmark with the same name moves it (removing the old instance).
$r = $d.Content
$r.Find.Execute("rival nationalist and ethnic factions in the FRY, which", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerEnd = $r.End
$bmRange = $d.Range($markerEnd, $markerEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
